$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(368).Insert()

$ws.Cells.Item(368, 1).Value = 5
$ws.Cells.Item(368, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(368, 3).Value = "Maule"
$ws.Cells.Item(368, 4).Value = 44951
$ws.Cells.Item(368, 5).Value = 7
$ws.Cells.Item(368, 6).Value = "Fruta"
$ws.Cells.Item(368, 7).Value = 100102
$ws.Cells.Item(368, 8).Value = "Cítricos"
$ws.Cells.Item(368, 9).Value = 100102004
$ws.Cells.Item(368, 10).Value = "Mandarina"
$ws.Cells.Item(368, 11).Value = "Murcott"
$ws.Cells.Item(368, 12).Value = "Primera"
$ws.Cells.Item(368, 13).Value = 230
$ws.Cells.Item(368, 14).Value = 9000
$ws.Cells.Item(368, 15).Value = 9000
$ws.Cells.Item(368, 16).Value = 9000
$ws.Cells.Item(368, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(368, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(368, 19).Value = 500
$ws.Cells.Item(368, 20).Value = 18
